$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2330.3333
$ws.Range("I53").Value = 3300
$ws.Range("J53").Value = 391
$ws.Range("K53").Value = 3300
$ws.Range("L53").Value = 391
$ws.Range("M53").Value = -2663
$ws.Range("N53").Value = -1665
$ws.Range("H74").Value = 3680
$ws.Range("H76").Value = 4640
$ws.Range("J76").Value = 5733.3335
$ws.Range("L76").Value = 5733.3335
$ws.Range("N76").Value = -6363.3335
$ws.Range("H77").Value = 3680
$ws.Range("H79").Value = 4640
$ws.Range("J79").Value = 5733.3335
$ws.Range("L79").Value = 5733.3335
$ws.Range("N79").Value = -7917.3335
$ws.Range("H98").Value = 3873.4167
$ws.Range("I98").Value = 4304.7095
$ws.Range("J98").Value = 1199.4
$ws.Range("K98").Value = 4304.7095
$ws.Range("L98").Value = 1199.4
$ws.Range("M98").Value = -2806.7095
$ws.Range("N98").Value = -4195.4
$ws.Range("H122").Value = 3873.4167
$ws.Range("I122").Value = 4304.7095
$ws.Range("J122").Value = 1199.4
$ws.Range("K122").Value = 12914.1285
$ws.Range("L122").Value = 3598.2
$ws.Range("M122").Value = -10464.1285
$ws.Range("N122").Value = -8498.200000000001
$ws.Range("H138").Value = 2614.4482
$ws.Range("I138").Value = 3024.625
$ws.Range("J138").Value = 2572.9114
$ws.Range("K138").Value = 9073.875
$ws.Range("L138").Value = 7718.7342
$ws.Range("M138").Value = -3933.875
$ws.Range("N138").Value = -17998.7342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 125002070
$ws.Range("I61").Value = 200001520
$ws.Range("K61").Value = 200001520
$ws.Range("M61").Value = -200001308
$ws.Range("H132").Value = 3750.4827
$ws.Range("I132").Value = 3470.5454
$ws.Range("K132").Value = 10411.6362
$ws.Range("M132").Value = -7881.636200000001
$ws.Range("H136").Value = 125002070
$ws.Range("I136").Value = 200001520
$ws.Range("K136").Value = 600004560
$ws.Range("M136").Value = -600002010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4252.4116
$ws.Range("J86").Value = 4126.75
$ws.Range("L86").Value = 4126.75
$ws.Range("N86").Value = -6372.75
$ws.Range("H89").Value = 4252.4116
$ws.Range("J89").Value = 4126.75
$ws.Range("L89").Value = 20633.75
$ws.Range("N89").Value = -31865.75
$ws.Range("H105").Value = 201982300
$ws.Range("I105").Value = 201982300
$ws.Range("K105").Value = 201982300
$ws.Range("M105").Value = -201980553

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 208.33333
$ws.Range("I7").Value = 93.40000000000001
$ws.Range("K7").Value = 93.40000000000001
$ws.Range("M7").Value = 19.59999999999999
$ws.Range("H22").Value = 63890
$ws.Range("I22").Value = 252
$ws.Range("J22").Value = 116921.664
$ws.Range("K22").Value = 252
$ws.Range("L22").Value = 116921.664
$ws.Range("M22").Value = 98
$ws.Range("N22").Value = -117621.664
$ws.Range("H58").Value = 11585.182
$ws.Range("I58").Value = 1980
$ws.Range("K58").Value = 1980
$ws.Range("M58").Value = -1777
$ws.Range("H134").Value = 16130497
$ws.Range("I134").Value = 1458.3914
$ws.Range("K134").Value = 4375.174199999999
$ws.Range("M134").Value = -1840.174199999999
$ws.Range("H136").Value = 11585.182
$ws.Range("I136").Value = 1980
$ws.Range("K136").Value = 5940
$ws.Range("M136").Value = -3390

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3878545.2
$ws.Range("I4").Value = 426500
$ws.Range("J4").Value = 4645666.5
$ws.Range("K4").Value = 1279500
$ws.Range("L4").Value = 13936999.5
$ws.Range("M4").Value = -1279388
$ws.Range("N4").Value = -13937223.5
$ws.Range("H5").Value = 361.93103
$ws.Range("I5").Value = 303.25
$ws.Range("K5").Value = 909.75
$ws.Range("M5").Value = -797.75
$ws.Range("H122").Value = 996.5
$ws.Range("J122").Value = 1178.0385
$ws.Range("L122").Value = 10602.3465
$ws.Range("N122").Value = -15502.3465
$ws.Range("H135").Value = 361.93103
$ws.Range("I135").Value = 303.25
$ws.Range("K135").Value = 2729.25
$ws.Range("M135").Value = -194.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56253220
$ws.Range("I70").Value = 35717964
$ws.Range("K70").Value = 35717964
$ws.Range("M70").Value = -35717694
$ws.Range("H73").Value = 56253220
$ws.Range("I73").Value = 35717964
$ws.Range("K73").Value = 35717964
$ws.Range("M73").Value = -35717028
$ws.Range("H80").Value = 4357
$ws.Range("I80").Value = 1900
$ws.Range("J80").Value = 4766.5
$ws.Range("K80").Value = 1900
$ws.Range("L80").Value = 4766.5
$ws.Range("M80").Value = -902
$ws.Range("N80").Value = -6762.5
$ws.Range("H83").Value = 4357
$ws.Range("I83").Value = 1900
$ws.Range("J83").Value = 4766.5
$ws.Range("K83").Value = 9500
$ws.Range("L83").Value = 23832.5
$ws.Range("M83").Value = -4508
$ws.Range("N83").Value = -33816.5
$ws.Range("H132").Value = 5288.3887
$ws.Range("I132").Value = 6548.3184
$ws.Range("K132").Value = 19644.9552
$ws.Range("M132").Value = -17114.9552
$ws.Range("H141").Value = 54299.875
$ws.Range("J141").Value = 54299.875
$ws.Range("L141").Value = 54299.875
$ws.Range("N141").Value = -64659.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 285.125
$ws.Range("I55").Value = 222.96
$ws.Range("K55").Value = 222.96
$ws.Range("M55").Value = -49.96000000000001
$ws.Range("H68").Value = 1270.5714
$ws.Range("I68").Value = 1252.9231
$ws.Range("K68").Value = 1252.9231
$ws.Range("M68").Value = -503.9231
$ws.Range("H71").Value = 1270.5714
$ws.Range("I71").Value = 1252.9231
$ws.Range("K71").Value = 6264.6155
$ws.Range("M71").Value = -2520.6155
$ws.Range("H100").Value = 1897.7778
$ws.Range("I100").Value = 1945
$ws.Range("J100").Value = 1860
$ws.Range("K100").Value = 1945
$ws.Range("L100").Value = 1860
$ws.Range("M100").Value = -1404
$ws.Range("N100").Value = -2942

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 83335630
$ws.Range("I62").Value = 250001500
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 250001500
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -250000876
$ws.Range("N62").Value = -3948
$ws.Range("H65").Value = 83335630
$ws.Range("I65").Value = 250001500
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 1250007500
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -1250004380
$ws.Range("N65").Value = -19740
$ws.Range("H141").Value = 57030
$ws.Range("J141").Value = 57030
$ws.Range("L141").Value = 57030
$ws.Range("N141").Value = -67390
